$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$win = $excel.ActiveWindow
$win.ScrollColumn = 12
$v = $win.ScrollColumn()
Write-Host "ScrollColumn value:" $v
$v2 = $win.ScrollRow()
Write-Host "ScrollRow value:" $v2
